# Updates workbook "OFICINA-CATAECSA" with 3 new advisor/client rows
# (ARQUITECKSA S.A., SOLIS SOLIS JUAN CARLOS, ZAVALA MENOSCAL HOMERO MIGUEL)
# across the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, plus a few
# related value updates (ARMIJO AGUILAR ROBERT LENIN and VEHINVER SA get
# additional sales; the "x de N" / total summary rows are refreshed).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"  (columns C:R)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# New row: ARQUITECKSA S.A. (inserted between ARMIJO and AVILA, i.e. row 7)
$ws1.Rows.Item(7).Insert()
$ws1.Range("A7").Value = "OFICINA-CATAECSA"
$ws1.Range("B7").Value = "ARQUITECKSA S.A."
$ws1.Range("C7:R7").Value = 0
$ws1.Range("D7").Value = 236.29

# New row: SOLIS SOLIS JUAN CARLOS (inserted between SOLIS OCAMPO and VEHINVER, row 24)
$ws1.Rows.Item(24).Insert()
$ws1.Range("A24").Value = "OFICINA-CATAECSA"
$ws1.Range("B24").Value = "SOLIS SOLIS JUAN CARLOS"
$ws1.Range("C24:R24").Value = 0
$ws1.Range("M24").Value = 3722.63

# New row: ZAVALA MENOSCAL HOMERO MIGUEL (inserted after VEHINVER, row 26)
$ws1.Rows.Item(26).Insert()
$ws1.Range("A26").Value = "OFICINA-CATAECSA"
$ws1.Range("B26").Value = "ZAVALA MENOSCAL HOMERO MIGUEL"
$ws1.Range("C26:R26").Value = 0

# ARMIJO AGUILAR ROBERT LENIN gained a PIEDRA SINTERIZADA sale
$ws1.Range("L6").Value = 1144.75

# VEHINVER SA gained a FREGADEROS DE COCINA + PORCELANATO sale
$ws1.Range("E25").Value = 94.20999999999999
$ws1.Range("M25").Value = 428.61

# Refresh the "x de 25" completion-count summary row (now row 27, 25 people)
$ws1.Range("C27").Value = "0 de 25"
$ws1.Range("D27").Value = "1 de 25"
$ws1.Range("E27").Value = "1 de 25"
$ws1.Range("F27").Value = "0 de 25"
$ws1.Range("G27").Value = "0 de 25"
$ws1.Range("H27").Value = "0 de 25"
$ws1.Range("I27").Value = "0 de 25"
$ws1.Range("J27").Value = "0 de 25"
$ws1.Range("K27").Value = "0 de 25"
$ws1.Range("L27").Value = "1 de 25"
$ws1.Range("M27").Value = "3 de 25"
$ws1.Range("N27").Value = "0 de 25"
$ws1.Range("O27").Value = "0 de 25"
$ws1.Range("P27").Value = "0 de 25"
$ws1.Range("Q27").Value = "0 de 25"
$ws1.Range("R27").Value = "0 de 25"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL"  (columns C:G)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# New row: ARQUITECKSA S.A.
$ws2.Rows.Item(7).Insert()
$ws2.Range("A7").Value = "OFICINA-CATAECSA"
$ws2.Range("B7").Value = "ARQUITECKSA S.A."
$ws2.Range("C7:G7").Value = 0
$ws2.Range("F7").Value = 236.29

# New row: SOLIS SOLIS JUAN CARLOS
$ws2.Rows.Item(24).Insert()
$ws2.Range("A24").Value = "OFICINA-CATAECSA"
$ws2.Range("B24").Value = "SOLIS SOLIS JUAN CARLOS"
$ws2.Range("C24:G24").Value = 0
$ws2.Range("F24").Value = 3722.63

# New row: ZAVALA MENOSCAL HOMERO MIGUEL
$ws2.Rows.Item(26).Insert()
$ws2.Range("A26").Value = "OFICINA-CATAECSA"
$ws2.Range("B26").Value = "ZAVALA MENOSCAL HOMERO MIGUEL"
$ws2.Range("C26:G26").Value = 0

# ARMIJO AGUILAR ROBERT LENIN: septiembre sale grew
$ws2.Range("F6").Value = 3842.16

# VEHINVER SA: septiembre sale added
$ws2.Range("F25").Value = 522.8200000000001

# Refresh the monthly totals row (now row 27)
$ws2.Range("C27").Value = 1076.87
$ws2.Range("D27").Value = 2057.76
$ws2.Range("E27").Value = 1423.94
$ws2.Range("F27").Value = 8323.9
$ws2.Range("G27").Value = 0
